$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8-16: item rows shift up by one slot, with several qty/rate/amount
# values updated in the process (per the source diff). ---

# Row 8
$ws.Range("A8").Formula = '="P. point"'
$ws.Range("A8").Copy()
$ws.Range("A8").PasteSpecial(-4163)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 78
$ws.Range("D8").Formula = '="2"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Formula = '="Short point (up to 3 mtr.)"'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)
$ws.Range("F8").Value = 256
$ws.Range("G8").Formula = '="19968.00"'
$ws.Range("G8").Copy()
$ws.Range("G8").PasteSpecial(-4163)
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "'"

# Row 9
$ws.Range("A9").Value = "'"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 56
$ws.Range("D9").Formula = '="3"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="Medium point (up to 6 mtr.)"'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)
$ws.Range("F9").Value = 472
$ws.Range("G9").Formula = '="26432.00"'
$ws.Range("G9").Copy()
$ws.Range("G9").PasteSpecial(-4163)
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = "'"

# Row 10
$ws.Range("A10").Formula = '="P. point"'
$ws.Range("A10").Copy()
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 74
$ws.Range("D10").Formula = '="4"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="Long point  (up to 10 mtr.)"'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)
$ws.Range("F10").Value = 662
$ws.Range("G10").Formula = '="48988.00"'
$ws.Range("G10").Copy()
$ws.Range("G10").PasteSpecial(-4163)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = "'"

# Row 11
$ws.Range("A11").Value = "'"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 34
$ws.Range("D11").Formula = '="2.0"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR"'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)
$ws.Range("F11").Value = 0
$ws.Range("G11").Formula = '="0.00"'
$ws.Range("G11").Copy()
$ws.Range("G11").PasteSpecial(-4163)
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "'"

# Row 12
$ws.Range("A12").Formula = '="P. point"'
$ws.Range("A12").Copy()
$ws.Range("A12").PasteSpecial(-4163)
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 52
$ws.Range("D12").Formula = '="6"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="On board"'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)
$ws.Range("F12").Value = 136
$ws.Range("G12").Formula = '="7072.00"'
$ws.Range("G12").Copy()
$ws.Range("G12").PasteSpecial(-4163)
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = "'"

# Row 13
$ws.Range("A13").Formula = '="Each"'
$ws.Range("A13").Copy()
$ws.Range("A13").PasteSpecial(-4163)
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 77
$ws.Range("D13").Formula = '="3.0"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR ."'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)
$ws.Range("F13").Value = 23
$ws.Range("G13").Formula = '="1771.00"'
$ws.Range("G13").Copy()
$ws.Range("G13").PasteSpecial(-4163)
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "'"

# Row 14
$ws.Range("A14").Value = "'"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 17
$ws.Range("D14").Formula = '="8"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Formula = '="Total"'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)
$ws.Range("F14").Value = 0
$ws.Range("G14").Formula = '="0.00"'
$ws.Range("G14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "'"

# Row 15
$ws.Range("A15").Formula = '="%"'
$ws.Range("A15").Copy()
$ws.Range("A15").PasteSpecial(-4163)
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 10
$ws.Range("D15").Formula = '="9"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="Add Tender Premium "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("F15").Value = 0
$ws.Range("G15").Formula = '="0.00"'
$ws.Range("G15").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = "'"

# Row 16
$ws.Range("A16").Value = "'"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 17
$ws.Range("D16").Formula = '="10"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Formula = '="Grand Total"'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("F16").Value = 0
$ws.Range("G16").Formula = '="0.00"'
$ws.Range("G16").Copy()
$ws.Range("G16").PasteSpecial(-4163)
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = "'"

# Row 17 keeps only its (blank) column A; B:I are cleared out entirely
$ws.Range("B17:I17").ClearContents()

# Row 18 becomes the new "Grand Total Rs." summary row
$ws.Range("A18").Value = "'"
$ws.Range("B18").Value = "'"
$ws.Range("C18").Value = "'"
$ws.Range("D18").Value = "'"
$ws.Range("E18").Formula = '="Grand Total Rs."'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("F18").Value = "'"
$ws.Range("G18").Formula = '="104231.00"'
$ws.Range("G18").Copy()
$ws.Range("G18").PasteSpecial(-4163)
$ws.Range("H18").Formula = '="104231.00"'
$ws.Range("H18").Copy()
$ws.Range("H18").PasteSpecial(-4163)
$ws.Range("I18").Value = "'"

# Row 19 becomes "Tender Premium @ 0%" (0.00 / 0.00)
$ws.Range("E19").Formula = '="Tender Premium @ 0%"'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)
$ws.Range("G19").Formula = '="0.00"'
$ws.Range("G19").Copy()
$ws.Range("G19").PasteSpecial(-4163)
$ws.Range("H19").Formula = '="0.00"'
$ws.Range("H19").Copy()
$ws.Range("H19").PasteSpecial(-4163)

# Row 20 becomes "NET PAYABLE AMOUNT Rs." (104231.00 / 104231.00)
$ws.Range("E20").Formula = '="NET PAYABLE AMOUNT Rs."'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)
$ws.Range("G20").Formula = '="104231.00"'
$ws.Range("G20").Copy()
$ws.Range("G20").PasteSpecial(-4163)
$ws.Range("H20").Formula = '="104231.00"'
$ws.Range("H20").Copy()
$ws.Range("H20").PasteSpecial(-4163)

# Old row 21 ("NET PAYABLE AMOUNT Rs." / 77354.00) is gone; the sheet now ends at row 20
$ws.Range("A21:I21").ClearContents()
